$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.981.63'
$ws.Range('E2').Value = '  +2.94%  '
$ws.Range('D3').Value = '2.537.15'
$ws.Range('E3').Value = '  +4.80%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''526.80'
$ws.Range('E5').Value = '  +2.74%  '
$ws.Range('D6').Value = '''134.99'
$ws.Range('E6').Value = '  +4.79%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').Value = '''0.567'
$ws.Range('E8').Value = '  +3.44%  '
$ws.Range('D9').Value = '2.534.04'
$ws.Range('D10').Value = '''0.0991'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').Value = '''5.21'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').Value = '2.984.75'
$ws.Range('E14').Value = '  +4.67%  '
$ws.Range('D15').Value = '58.949.17'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').Value = '''22.35'
$ws.Range('E16').Value = '  +4.16%  '
$ws.Range('E17').Value = '  +3.73%  '
$ws.Range('D18').Value = '2.535.04'
$ws.Range('E18').Value = '  +4.50%  '
$ws.Range('D19').Value = '''10.75'
$ws.Range('E19').Value = '  +3.72%  '
$ws.Range('D20').Value = '''323.97'
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('D21').Value = '''4.20'
$ws.Range('E21').Value = '  +3.21%  '
$ws.Range('D22').Value = '''6.11'
$ws.Range('E22').Value = '  +8.52%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '''65.14'
$ws.Range('E24').Value = '  +2.59%  '
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('E28').Value = '  +4.44%  '
$ws.Range('E29').Value = '  +6.21%  '
$ws.Range('E30').Value = '  +8.35%  '
$ws.Range('E31').Value = '  +4.84%  '
$ws.Range('D32').Value = '''169.98'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '''6.35'
$ws.Range('E33').Value = '  +2.74%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '''0.997'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '''18.28'
$ws.Range('E36').Value = '  +3.58%  '
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('E39').Value = '  +5.20%  '
$ws.Range('D40').Value = '''36.75'
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('D41').Value = '''0.786'
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('D42').Value = '''280.17'
$ws.Range('E42').Value = '  +6.16%  '
$ws.Range('D43').Value = '''134.85'
$ws.Range('E43').Value = '  +11.68%  '
$ws.Range('E44').Value = '  +3.69%  '
$ws.Range('D45').Value = '''5.09'
$ws.Range('E45').Value = '  +5.27%  '
$ws.Range('D46').Value = '''0.602'
$ws.Range('E46').Value = '  +3.47%  '
$ws.Range('D47').Value = '''0.0923'
$ws.Range('E47').Value = '  +2.53%  '
$ws.Range('E49').Value = '  +4.05%  '
$ws.Range('E50').Value = '  +4.36%  '
$ws.Range('D51').Value = '1.757.51'
$ws.Range('E51').Value = '  +4.11%  '
